$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-3 and shift/replace row4 onward with new data,
# extending the table from A1:D4 to A1:D6.

$data = @(
    @(0.25, 88924.074024154499, 61510.287898985996, 179.96238143455),
    @(0.5,  87646.622636007596, 60264.573801184903, 211.69967177987201),
    @(1,    85104.733203326701, 57788.213950614001, 277.22925389091398),
    @(2,    80744.775875134001, 53562.463715507598, 411.43634697693301),
    @(4,    75683.053847636693, 48747.494175761502, 658.18883472863001)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("D21").Select()
